$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4587.052556530442
$ws.Range("D3").Value = 4587.052556530438
$ws.Range("D4").Value = 4587.052556530438
$ws.Range("D6").Value = 13105.86444722955
$ws.Range("D7").Value = 13105.86444722955
$ws.Range("D9").Value = 1269.52664843566
$ws.Range("D10").Value = 1269.526648435659
$ws.Range("D11").Value = 18.46262876487648
$ws.Range("D12").Value = 1047.975103257141
$ws.Range("D13").Value = 18.46262876487648
$ws.Range("D17").Value = 14770.10301190118
$ws.Range("D19").Value = 11888.26840371751
$ws.Range("D20").Value = 11888.26840371751
$ws.Range("D21").Value = 221.5515451785177
$ws.Range("D24").Value = 207268.5707538992
$ws.Range("D25").Value = 207268.570753899
$ws.Range("D28").Value = 4145.371415078048
$ws.Range("D29").Value = 4145.371415078048
$ws.Range("D30").Value = 207268.570753899
$ws.Range("D35").Value = 20217.60000000001
$ws.Range("D36").Value = 20217.60000000001
$ws.Range("D38").Value = -237.5265709177201
$ws.Range("D39").Value = -237.5265709177199
$ws.Range("D41").Value = 237.5265709177199
$ws.Range("D42").Value = 4750.531418354351
$ws.Range("D43").Value = 4750.531418354351
